# Fruta / hortaliza, semanal
# This weekly refresh re-sorts the Higo (fig) price rows by their actual
# market date. The row data for columns D (Fecha), M (Volumen),
# N (Precio mínimo), O (Precio máximo), P (Precio promedio ponderado),
# R (Origen) and S (Precio $/Kg) gets rotated between rows while the rest
# of each row (Mercado, Producto, Calidad, Unidad de comercialización, ...)
# stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as a "data block" per source row.
$cols = @("D", "M", "N", "O", "P", "R", "S")

# Snapshot the current (pre-edit) values for every row that participates
# in the re-sort.
$rows = @(2, 3, 4, 5, 6, 7, 10, 11, 12, 13, 14, 15)
$snapshot = @{}
foreach ($r in $rows) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range($c + $r).Value()
    }
    $snapshot[$r] = $vals
}

# target row -> source row (where the "new" data for that row comes from)
$mapping = @{
    2  = 10
    3  = 11
    4  = 14
    5  = 15
    6  = 4
    7  = 5
    10 = 2
    11 = 3
    12 = 6
    13 = 7
    14 = 12
    15 = 13
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $srcVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range($c + $targetRow).Value = $srcVals[$c]
    }
}
